$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '68.011.68'
$ws.Range("E2").Value = '  +1.07%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '2.628.35'
$ws.Range("E3").Value = '  +0.13%  '

# Row 4: update D4, E4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5: update D5, E5
$ws.Range("D5").Value = '''598.32'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6: update D6, E6
$ws.Range("D6").Value = '''153.15'
$ws.Range("E6").Value = '  +0.39%  '

# Row 7: update E7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8: update E8
$ws.Range("E8").Value = '  -1.87%  '

# Row 9: update D9, E9
$ws.Range("D9").Value = '2.625.94'
$ws.Range("E9").Value = '  +0.12%  '

# Row 10: update E10
$ws.Range("E10").Value = '  +9.14%  '

# Row 11: update E11
$ws.Range("E11").Value = '  -0.51%  '

# Row 12: update D12, E12
$ws.Range("D12").Value = '''5.21'
$ws.Range("E12").Value = '  +0.23%  '

# Row 13: update D13, E13
$ws.Range("D13").Value = '''0.347'
$ws.Range("E13").Value = '  -0.86%  '

# Row 14: update D14, E14
$ws.Range("D14").Value = '''27.65'
$ws.Range("E14").Value = '  +0.35%  '

# Row 15: update E15
$ws.Range("E15").Value = '  +3.70%  '

# Row 16: update D16, E16
$ws.Range("D16").Value = '3.105.00'
$ws.Range("E16").Value = '  +0.09%  '

# Row 17: update D17, E17
$ws.Range("D17").Value = '67.831.55'
$ws.Range("E17").Value = '  +0.87%  '

# Row 18: update D18, E18
$ws.Range("D18").Value = '2.630.95'
$ws.Range("E18").Value = '  +0.25%  '

# Row 19: update D19, E19
$ws.Range("D19").Value = '''372.92'
$ws.Range("E19").Value = '  +2.53%  '

# Row 20: update D20, E20
$ws.Range("D20").Value = '''11.29'
$ws.Range("E20").Value = '  +1.29%  '

# Row 21: update D21, E21
$ws.Range("D21").Value = '''7.48'
$ws.Range("E21").Value = '  -0.19%  '

# Row 22: update D22, E22
$ws.Range("D22").Value = '''4.24'
$ws.Range("E22").Value = '  -1.21%  '

# Row 23: update D23, E23
$ws.Range("D23").Value = '''4.82'
$ws.Range("E23").Value = '  -1.83%  '

# Row 24: update E24
$ws.Range("E24").Value = '  -2.72%  '

# Row 25: update D25, E25
$ws.Range("D25").Value = '''72.69'
$ws.Range("E25").Value = '  +9.64%  '

# Row 26: update E26
$ws.Range("E26").Value = '  +0.08%  '

# Row 27: update D27, E27
$ws.Range("D27").Value = '''9.91'
$ws.Range("E27").Value = '  -2.12%  '

# Row 28: update D28, E28
$ws.Range("D28").Value = '''0.0000105'
$ws.Range("E28").Value = '  +2.54%  '

# Row 29: update D29, E29
$ws.Range("D29").Value = '2.756.60'
$ws.Range("E29").Value = '  -0.26%  '

# Row 30: update D30, E30
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -2.23%  '

# Row 31: update D31, E31
$ws.Range("D31").Value = '''579.27'
$ws.Range("E31").Value = '  +0.51%  '

# Row 32: update D32, E32
$ws.Range("D32").Value = '''1.40'
$ws.Range("E32").Value = '  +0.56%  '

# Row 33: update D33, E33
$ws.Range("D33").Value = '''7.82'
$ws.Range("E33").Value = '  +0.48%  '

# Row 34: update D34, E34
$ws.Range("D34").Value = '''1.85'
$ws.Range("E34").Value = '  +0.39%  '

# Row 35: update E35
$ws.Range("E35").Value = '  +0.02%  '

# Row 36: update D36, E36
$ws.Range("D36").Value = '''0.127'
$ws.Range("E36").Value = '  -0.90%  '

# Row 37: update E37
$ws.Range("E37").Value = '  -0.47%  '

# Row 38: update D38, E38
$ws.Range("D38").Value = '''158.99'
$ws.Range("E38").Value = '  +0.74%  '

# Row 39: update E39
$ws.Range("E39").Value = '  -0.49%  '

# Row 40: update D40, E40
$ws.Range("D40").Value = '''1.90'
$ws.Range("E40").Value = '  +4.54%  '

# Row 41: update D41, E41
$ws.Range("D41").Value = '''0.369'
$ws.Range("E41").Value = '  +0.41%  '

# Row 42: update D42, E42
$ws.Range("D42").Value = '''5.33'
$ws.Range("E42").Value = '  +1.37%  '

# Row 43: update D43, E43
$ws.Range("D43").Value = '''2.63'
$ws.Range("E43").Value = '  +1.49%  '

# Row 44: update D44, E44
$ws.Range("D44").Value = '''17.10'
$ws.Range("E44").Value = '  +4.58%  '

# Row 45: update D45, E45
$ws.Range("D45").Value = '0.0₆0316'
$ws.Range("E45").Value = '  +11.10%  '

# Row 46: update D46, E46
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  +0.10%  '

# Row 47: update D47, E47
$ws.Range("D47").Value = '''40.42'
$ws.Range("E47").Value = '  -1.88%  '

# Row 48: update D48, E48
$ws.Range("D48").Value = '''155.56'
$ws.Range("E48").Value = '  +0.29%  '

# Row 49: update D49, E49
$ws.Range("D49").Value = '''3.70'
$ws.Range("E49").Value = '  -0.69%  '

# Row 50: update D50, E50
$ws.Range("D50").Value = '''1.70'
$ws.Range("E50").Value = '  -0.97%  '

# Row 51: update B51, C51, D51, E51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''21.80'
$ws.Range("E51").Value = '  +6.70%  '
